# Weekly update: insert a new week of Tomate price rows for
# "Comercializadora del Agro de Limarí" ahead of the existing history,
# shifting the rest of the table down by 6 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new blank rows at the top of the price block (rows 471-476),
# pushing the existing rows 471-501 down to 477-507.
$ws.Range("A471:R476").EntireRow.Insert()

# New week's data (date serial 44615 = 2022-02-23).
$rows = @(
    @{ Row = 471; H = "Larga vida"; I = "Primera"; J = 2400; K = 9000;  L = 10000; M = 9500; P = 528 },
    @{ Row = 472; H = "Larga vida"; I = "Segunda"; J = 1900; K = 7000;  L = 8000;  M = 7500; P = 417 },
    @{ Row = 473; H = "Larga vida"; I = "Tercera"; J = 600;  K = 5000;  L = 6000;  M = 5500; P = 306 },
    @{ Row = 474; H = "Semiduro";   I = "Primera"; J = 2400; K = 6000;  L = 7000;  M = 6500; P = 361 },
    @{ Row = 475; H = "Semiduro";   I = "Segunda"; J = 1800; K = 4000;  L = 5000;  M = 4500; P = 250 },
    @{ Row = 476; H = "Semiduro";   I = "Tercera"; J = 1300; K = 2000;  L = 3000;  M = 2500; P = 139 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 1).Value = 2
    $ws.Cells.Item($n, 2).Value = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($n, 3).Value = "Coquimbo"
    $ws.Cells.Item($n, 4).Value = 44615
    $ws.Cells.Item($n, 5).Value = 4
    $ws.Cells.Item($n, 6).Value = 100112020
    $ws.Cells.Item($n, 7).Value = "Tomate"
    $ws.Cells.Item($n, 8).Value = $r.H
    $ws.Cells.Item($n, 9).Value = $r.I
    $ws.Cells.Item($n, 10).Value = $r.J
    $ws.Cells.Item($n, 11).Value = $r.K
    $ws.Cells.Item($n, 12).Value = $r.L
    $ws.Cells.Item($n, 13).Value = $r.M
    $ws.Cells.Item($n, 14).Value = "$/bandeja 18 kilos"
    $ws.Cells.Item($n, 15).Value = "Provincia de Limarí"
    $ws.Cells.Item($n, 16).Value = $r.P
    $ws.Cells.Item($n, 17).Value = 18
    $ws.Cells.Item($n, 18).Value = "Hortaliza"
}
